# Insert one new weekly record at row 34 (Macroferia Regional de Talca - Choclo),
# which pushes the existing rows 34-130 down to 35-131.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(34).Insert()

$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(34, 3).Value = "Maule"
$ws.Cells.Item(34, 4).Value = 44536
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = 100112024
$ws.Cells.Item(34, 7).Value = "Choclo"
$ws.Cells.Item(34, 8).Value = "Choclero"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 20000
$ws.Cells.Item(34, 11).Value = 400
$ws.Cells.Item(34, 12).Value = 400
$ws.Cells.Item(34, 13).Value = 400
$ws.Cells.Item(34, 14).Value = "$/unidad"
$ws.Cells.Item(34, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(34, 16).Value = 400
$ws.Cells.Item(34, 17).Value = 1
$ws.Cells.Item(34, 18).Value = "Hortaliza"
